$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = -0.3125041300915371
$ws.Range("J4").Value = 0.4712093240364904
$ws.Range("K4").Value = 0.7536979182829137
$ws.Range("L4").Value = 3.332998588019548
